$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 27 (cohort 2021, period 4): num_customers 58 -> 59, retention_rate recalculated
$ws.Range("C27").Value = 59
$ws.Range("E27").Value = 59/2252

# Row 34 (cohort 2022, period 2): num_customers 91 -> 93, retention_rate recalculated
$ws.Range("C34").Value = 93
$ws.Range("E34").Value = 93/2256

# Row 36 (cohort 2023, period 1): num_customers 149 -> 152, retention_rate recalculated
$ws.Range("C36").Value = 152
$ws.Range("E36").Value = 152/1930

# Row 37 (cohort 2024, period 0): num_customers and cohort_size 1014 -> 1025
$ws.Range("C37").Value = 1025
$ws.Range("D37").Value = 1025
